# Auto-generated script applying updated market-price figures to the Ixion Profits workbook
# (columns H-N on various rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3731222
$ws.Range("J17").Value = 3731222
$ws.Range("L17").Value = 11193666
$ws.Range("N17").Value = -11194002

$ws.Range("H106").Value = 47620800
$ws.Range("I106").Value = 22223726
$ws.Range("K106").Value = 22223726
$ws.Range("M106").Value = -22223095

$ws.Range("H116").Value = 9199.333000000001
$ws.Range("I116").Value = 34568.332
$ws.Range("J116").Value = 2857.0833
$ws.Range("K116").Value = 34568.332
$ws.Range("L116").Value = 2857.0833
$ws.Range("M116").Value = -31126.332
$ws.Range("N116").Value = -9741.0833

$ws.Range("H123").Value = 30641.111
$ws.Range("J123").Value = 30641.111
$ws.Range("L123").Value = 30641.111
$ws.Range("N123").Value = -40441.111

$ws.Range("H125").Value = 3973.158
$ws.Range("I125").Value = 12376
$ws.Range("J125").Value = 2984.5881
$ws.Range("K125").Value = 111384
$ws.Range("L125").Value = 26861.2929
$ws.Range("M125").Value = -108924
$ws.Range("N125").Value = -31781.2929

$ws.Range("H137").Value = 1279.4783
$ws.Range("I137").Value = 1171.0278
$ws.Range("J137").Value = 1669.9
$ws.Range("K137").Value = 3513.0834
$ws.Range("L137").Value = 5009.700000000001
$ws.Range("M137").Value = -963.0834000000004
$ws.Range("N137").Value = -10109.7

$ws.Range("H138").Value = 3485.191
$ws.Range("I138").Value = 795.0577
$ws.Range("J138").Value = 7265.919
$ws.Range("K138").Value = 2385.1731
$ws.Range("L138").Value = 21797.757
$ws.Range("M138").Value = 2754.8269
$ws.Range("N138").Value = -32077.757

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5675.76
$ws.Range("I32").Value = 4070.9395
$ws.Range("K32").Value = 4070.9395
$ws.Range("M32").Value = -3783.9395

$ws.Range("H61").Value = 3363.98
$ws.Range("I61").Value = 3295.8125
$ws.Range("K61").Value = 3295.8125
$ws.Range("M61").Value = -3083.8125

$ws.Range("H136").Value = 3363.98
$ws.Range("I136").Value = 3295.8125
$ws.Range("K136").Value = 9887.4375
$ws.Range("M136").Value = -7337.4375

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 65450.5
$ws.Range("J58").Value = 65450.5
$ws.Range("L58").Value = 65450.5
$ws.Range("N58").Value = -66038.5

$ws.Range("H59").Value = 47890
$ws.Range("J59").Value = 47890
$ws.Range("L59").Value = 47890
$ws.Range("N59").Value = -49584

$ws.Range("H60").Value = 32500
$ws.Range("J60").Value = 32500
$ws.Range("L60").Value = 32500
$ws.Range("N60").Value = -33698

$ws.Range("H94").Value = 2266.611
$ws.Range("I94").Value = 1889.9
$ws.Range("K94").Value = 1889.9
$ws.Range("M94").Value = -1438.9

$ws.Range("H134").Value = 3500.0527
$ws.Range("I134").Value = 3721.7556
$ws.Range("J134").Value = 2668.6667
$ws.Range("K134").Value = 11165.2668
$ws.Range("L134").Value = 8006.000100000001
$ws.Range("M134").Value = -8630.266799999999
$ws.Range("N134").Value = -13076.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 565.5
$ws.Range("I22").Value = 547.75
$ws.Range("J22").Value = 601
$ws.Range("K22").Value = 547.75
$ws.Range("L22").Value = 601
$ws.Range("M22").Value = -197.75
$ws.Range("N22").Value = -1301

$ws.Range("H31").Value = 13160864
$ws.Range("I31").Value = 1329.8928
$ws.Range("J31").Value = 50007556
$ws.Range("K31").Value = 1329.8928
$ws.Range("L31").Value = 50007556
$ws.Range("M31").Value = -1034.8928
$ws.Range("N31").Value = -50008146

$ws.Range("H34").Value = 13160864
$ws.Range("I34").Value = 1329.8928
$ws.Range("J34").Value = 50007556
$ws.Range("K34").Value = 1329.8928
$ws.Range("L34").Value = 50007556
$ws.Range("M34").Value = -1127.8928
$ws.Range("N34").Value = -50007960

$ws.Range("H58").Value = 5556868
$ws.Range("I58").Value = 6945395.5
$ws.Range("J58").Value = 2757
$ws.Range("K58").Value = 6945395.5
$ws.Range("L58").Value = 2757
$ws.Range("M58").Value = -6945192.5
$ws.Range("N58").Value = -3163

$ws.Range("H124").Value = 29999.75
$ws.Range("J124").Value = 29999.75
$ws.Range("L124").Value = 29999.75
$ws.Range("N124").Value = -34909.75

$ws.Range("H136").Value = 5556868
$ws.Range("I136").Value = 6945395.5
$ws.Range("J136").Value = 2757
$ws.Range("K136").Value = 20836186.5
$ws.Range("L136").Value = 8271
$ws.Range("M136").Value = -20833636.5
$ws.Range("N136").Value = -13371

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5101207
$ws.Range("I5").Value = 480.8
$ws.Range("J5").Value = 6922895
$ws.Range("K5").Value = 1442.4
$ws.Range("L5").Value = 20768685
$ws.Range("M5").Value = -1330.4
$ws.Range("N5").Value = -20768909

$ws.Range("H23").Value = 7142946
$ws.Range("I23").Value = 16666742
$ws.Range("J23").Value = 98.25
$ws.Range("K23").Value = 50000226
$ws.Range("L23").Value = 294.75
$ws.Range("M23").Value = -49999991
$ws.Range("N23").Value = -764.75

$ws.Range("H113").Value = 4286201
$ws.Range("I113").Value = 7143302.5
$ws.Range("K113").Value = 21429907.5
$ws.Range("M113").Value = -21427737.5

$ws.Range("H117").Value = 30308360
$ws.Range("I117").Value = 409.66666
$ws.Range("J117").Value = 41673840
$ws.Range("K117").Value = 1228.99998
$ws.Range("L117").Value = 125021520
$ws.Range("M117").Value = 2213.00002
$ws.Range("N117").Value = -125028404

$ws.Range("H132").Value = 6945694
$ws.Range("I132").Value = 974.8333
$ws.Range("K132").Value = 8773.4997
$ws.Range("M132").Value = -6243.4997

$ws.Range("H135").Value = 5101207
$ws.Range("I135").Value = 480.8
$ws.Range("J135").Value = 6922895
$ws.Range("K135").Value = 4327.2
$ws.Range("L135").Value = 62306055
$ws.Range("M135").Value = -1792.2
$ws.Range("N135").Value = -62311125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3481490.8
$ws.Range("I122").Value = 3088258.5
$ws.Range("K122").Value = 9264775.5
$ws.Range("M122").Value = -9262325.5

$ws.Range("H123").Value = 27799.6
$ws.Range("J123").Value = 27799.6
$ws.Range("L123").Value = 27799.6
$ws.Range("N123").Value = -32699.6

$ws.Range("H132").Value = 4388209
$ws.Range("I132").Value = 5378224
$ws.Range("K132").Value = 16134672
$ws.Range("M132").Value = -16132142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8743074
$ws.Range("I22").Value = 43704036
$ws.Range("J22").Value = 2832.6667
$ws.Range("K22").Value = 43704036
$ws.Range("L22").Value = 2832.6667
$ws.Range("M22").Value = -43703741
$ws.Range("N22").Value = -3422.6667

$ws.Range("H27").Value = 8743074
$ws.Range("I27").Value = 43704036
$ws.Range("J27").Value = 2832.6667
$ws.Range("K27").Value = 43704036
$ws.Range("L27").Value = 2832.6667
$ws.Range("M27").Value = -43703929
$ws.Range("N27").Value = -3046.6667

$ws.Range("H82").Value = 207660
$ws.Range("I82").Value = 3100
$ws.Range("J82").Value = 514500
$ws.Range("K82").Value = 3100
$ws.Range("L82").Value = 514500
$ws.Range("M82").Value = -2739
$ws.Range("N82").Value = -515222

$ws.Range("H85").Value = 207660
$ws.Range("I85").Value = 3100
$ws.Range("J85").Value = 514500
$ws.Range("K85").Value = 3100
$ws.Range("L85").Value = 514500
$ws.Range("M85").Value = -1852
$ws.Range("N85").Value = -516996

$ws.Range("H132").Value = 7642671.5
$ws.Range("I132").Value = 10917016
$ws.Range("J132").Value = 2533.389
$ws.Range("K132").Value = 32751048
$ws.Range("L132").Value = 7600.167
$ws.Range("M132").Value = -32748518
$ws.Range("N132").Value = -12660.167

$ws.Range("H136").Value = 8565.344999999999
$ws.Range("I136").Value = 7884.3
$ws.Range("J136").Value = 10078.777
$ws.Range("K136").Value = 23652.9
$ws.Range("L136").Value = 30236.331
$ws.Range("M136").Value = -21102.9
$ws.Range("N136").Value = -35336.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1013.21313
$ws.Range("I132").Value = 683.8163500000001
$ws.Range("J132").Value = 2358.25
$ws.Range("K132").Value = 2358.25
$ws.Range("L132").Value = 7074.75
$ws.Range("M132").Value = 478.5509499999998
$ws.Range("N132").Value = -12134.75
